$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the day-58 (four-square-ping station) row values that were
# previously blank.
$ws.Range("C58").Value = 1
$ws.Range("D58").Value = 1
$ws.Range("E58").Value = 1
$ws.Range("F58").Value = 1

# F58 had no prior number format (default style); give it the same
# "0.00_ " custom format used by the rest of column C:E so it picks up
# the shared style index instead of the default one.
$ws.Range("F58").NumberFormat = "0.00_ "

# Move the active selection to H56 (was I59).
$ws.Range("H56").Select()
